$d = $word.ActiveDocument

$replacements = @(
    @("310÷8=38, 6", "289÷2=144, 1"),
    @("310÷2=155, 0", "110÷6=18, 2"),
    @("914÷2=457, 0", "838÷7=119, 5"),
    @("820÷2=410, 0", "721÷9=80, 1"),
    @("642÷9=71, 3", "464÷8=58, 0"),
    @("514÷9=57, 1", "742÷6=123, 4"),
    @("570÷2=285, 0", "967÷3=322, 1"),
    @("999÷8=124, 7", "793÷9=88, 1"),
    @("584÷3=194, 2", "627÷3=209, 0"),
    @("577÷2=288, 1", "145÷6=24, 1"),
    @("107÷3=35, 2", "981÷8=122, 5"),
    @("397÷7=56, 5", "200÷9=22, 2"),
    @("606÷7=86, 4", "992÷9=110, 2"),
    @("247÷9=27, 4", "612÷7=87, 3"),
    @("484÷7=69, 1", "684÷5=136, 4"),
    @("603÷2=301, 1", "244÷4=61, 0"),
    @("648÷5=129, 3", "468÷5=93, 3"),
    @("642÷7=91, 5", "914÷7=130, 4"),
    @("166÷9=18, 4", "812÷3=270, 2"),
    @("345÷9=38, 3", "101÷4=25, 1"),
    @("553÷7=79, 0", "960÷3=320, 0"),
    @("981÷3=327, 0", "192÷4=48, 0"),
    @("445÷3=148, 1", "419÷7=59, 6"),
    @("175÷7=25, 0", "788÷8=98, 4"),
    @("388÷3=129, 1", "468÷9=52, 0")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
